# Updates the Price (D) / Volume(1h) (E) columns of the cryptos sheet with a
# fresh snapshot of quotes. Both columns store plain text in the workbook
# (prices like "523.17" are NOT numbers - they're already-formatted strings),
# so for D-column values that look like a plain number we prefix them with a
# leading apostrophe; that is the standard Excel "store as text" trick and
# keeps Range.Value from auto-coercing the literal into a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.707.25'
$ws.Range("E2").Value = '  +0.65%  '

$ws.Range("D3").Value = '3.102.53'
$ws.Range("E3").Value = '  +1.37%  '

$ws.Range("D5").Value = '''522.95'
$ws.Range("E5").Value = '  +1.36%  '

$ws.Range("D6").Value = '''140.45'
$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '3.102.46'
$ws.Range("E8").Value = '  +1.52%  '

$ws.Range("D9").Value = '''0.435'
$ws.Range("E9").Value = '  +0.72%  '

$ws.Range("E10").Value = '  -0.01%  '

$ws.Range("E11").Value = '  +1.05%  '

$ws.Range("E12").Value = '  +2.67%  '

$ws.Range("D13").Value = '3.636.16'
$ws.Range("E13").Value = '  +1.26%  '

$ws.Range("E14").Value = '  +1.65%  '

$ws.Range("D15").Value = '''26.16'
$ws.Range("E15").Value = '  +2.67%  '

$ws.Range("E16").Value = '  +0.92%  '

$ws.Range("D17").Value = '57.796.13'
$ws.Range("E17").Value = '  +0.68%  '

$ws.Range("D18").Value = '3.109.60'
$ws.Range("E18").Value = '  +1.16%  '

$ws.Range("E19").Value = '  +0.53%  '

$ws.Range("D20").Value = '''12.79'
$ws.Range("E20").Value = '  -1.09%  '

$ws.Range("D21").Value = '''8.01'
$ws.Range("E21").Value = '  -0.77%  '

$ws.Range("D22").Value = '''336.81'
$ws.Range("E22").Value = '  +1.25%  '

$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("E24").Value = '  +1.47%  '

$ws.Range("D25").Value = '''66.44'
$ws.Range("E25").Value = '  +1.01%  '

$ws.Range("E26").Value = '  -1.04%  '

$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("D28").Value = '0.0₃0925'
$ws.Range("E28").Value = '  +1.91%  '

$ws.Range("D29").Value = '''6.54'
$ws.Range("E29").Value = '  +3.16%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").Value = '''7.21'
$ws.Range("E31").Value = '  +0.65%  '

$ws.Range("E32").Value = '  +2.32%  '

$ws.Range("E33").Value = '  +0.65%  '

$ws.Range("D34").Value = '''1.19'
$ws.Range("E34").Value = '  +1.60%  '

$ws.Range("D35").Value = '''154.23'
$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("E36").Value = '  +3.83%  '

$ws.Range("D37").Value = '''6.07'
$ws.Range("E37").Value = '  +3.29%  '

$ws.Range("D38").Value = '''26.91'
$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("E39").Value = '  +2.34%  '

$ws.Range("D40").Value = '''0.0666'
$ws.Range("E40").Value = '  -0.80%  '

$ws.Range("D41").Value = '3.146.72'
$ws.Range("E41").Value = '  +1.50%  '

$ws.Range("D42").Value = '''0.683'
$ws.Range("E42").Value = '  +4.34%  '

$ws.Range("D43").Value = '''36.86'
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("E44").Value = '  -0.46%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("E46").Value = '  +5.96%  '

$ws.Range("D47").Value = '2.278.97'
$ws.Range("E47").Value = '  +0.76%  '

$ws.Range("E48").Value = '  +0.78%  '

$ws.Range("E49").Value = '  +6.36%  '

$ws.Range("D50").Value = '''20.61'
$ws.Range("E50").Value = '  +1.38%  '

$ws.Range("D51").Value = '''5.99'
$ws.Range("E51").Value = '  +2.47%  '

